$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'297.88"
$ws.Range("E2").Value = "'-4.68%"
$ws.Range("D3").Value = "'31.72"
$ws.Range("E3").Value = "'-1.09%"
$ws.Range("D4").Value = "'5.106"
$ws.Range("E4").Value = "'-4.31%"
$ws.Range("D5").Value = "'0.07528"
$ws.Range("E5").Value = "'-0.36%"
$ws.Range("D6").Value = "'7.753"
$ws.Range("E6").Value = "'-1.12%"
$ws.Range("D7").Value = "'1.716"
$ws.Range("E7").Value = "'9.15%"
$ws.Range("D8").Value = "'3.794"
$ws.Range("E8").Value = "'3.31%"
$ws.Range("D9").Value = "'0.9301"
$ws.Range("E9").Value = "'2.03%"
$ws.Range("E10").Value = "'0.04%"
$ws.Range("D11").Value = "'0.07492"
$ws.Range("E11").Value = "'-3.18%"
$ws.Range("D12").Value = "'0.08004"
$ws.Range("E12").Value = "'-1.44%"
$ws.Range("D13").Value = "'0.03032"
$ws.Range("E13").Value = "'0.07%"
$ws.Range("D14").Value = "'0.09895"
$ws.Range("E14").Value = "'0.26%"
$ws.Range("D15").Value = "'0.001507"
$ws.Range("E15").Value = "'-1.27%"
$ws.Range("D16").Value = "'0.006352"
$ws.Range("E16").Value = "'-0.76%"
$ws.Range("D17").Value = "'3.458"
$ws.Range("E17").Value = "'-1.26%"
$ws.Range("D18").Value = "'2.225"
$ws.Range("E18").Value = "'-0.71%"
$ws.Range("E19").Value = "'0.18%"
$ws.Range("E20").Value = "'0.80%"
$ws.Range("D21").Value = "'4.562"
$ws.Range("E21").Value = "'9.19%"
$ws.Range("D22").Value = "'0.04657"
$ws.Range("E22").Value = "'2.26%"
$ws.Range("D24").Value = "'0.001221"
$ws.Range("E24").Value = "'0.44%"
$ws.Range("D25").Value = "'0.004422"
$ws.Range("E25").Value = "'-1.53%"
$ws.Range("D26").Value = "'0.0001301"
$ws.Range("E26").Value = "'0.20%"
$ws.Range("E27").Value = "'6.98%"
$ws.Range("D39").Value = "'0.01671"
$ws.Range("E39").Value = "'-1.82%"
$ws.Range("D40").Value = "'0.04526"
$ws.Range("E40").Value = "'-1.06%"
$ws.Range("D41").Value = "'0.007101"
$ws.Range("E41").Value = "'-1.41%"
$ws.Range("D42").Value = "'0.1327"
$ws.Range("E42").Value = "'-2.64%"
$ws.Range("D43").Value = "'0.002061"
$ws.Range("E43").Value = "'-8.67%"
$ws.Range("D44").Value = "'0.01252"
$ws.Range("E44").Value = "'-10.07%"
$ws.Range("D45").Value = "'0.00005995"
$ws.Range("E45").Value = "'-2.86%"
$ws.Range("D46").Value = "'1.930"
$ws.Range("E46").Value = "'1.96%"
$ws.Range("E47").Value = "'-0.01%"
